# Auto-applied market-price refresh for Odin_Profits sheets.
# For each affected leve row, currentAveragePrice* / LevePrice* / LeveProfit*
# columns (H:N) are refreshed from the latest market snapshot.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 53
$ws.Range("H53").Value2 = 361.2143
$ws.Range("I53").Value2 = 104.5
$ws.Range("J53").Value2 = 1003
$ws.Range("K53").Value2 = 104.5
$ws.Range("L53").Value2 = 1003
$ws.Range("M53").Value2 = 532.5
$ws.Range("N53").Value2 = -2277

# Row 98
$ws.Range("H98").Value2 = 1907.75
$ws.Range("I98").Value2 = 1901.6666
$ws.Range("J98").Value2 = 1999
$ws.Range("K98").Value2 = 1901.6666
$ws.Range("L98").Value2 = 1999
$ws.Range("M98").Value2 = -403.6666
$ws.Range("N98").Value2 = -4995

# Row 101
$ws.Range("H101").Value2 = 2234
$ws.Range("I101").Value2 = 2234
$ws.Range("J101").Value2 = 0
$ws.Range("K101").Value2 = 6702
$ws.Range("L101").Value2 = 0
$ws.Range("M101").Value2 = -5080
$ws.Range("N101").ClearContents()

# Row 107
$ws.Range("H107").Value2 = 1845.7368
$ws.Range("J107").Value2 = 4183.8
$ws.Range("L107").Value2 = 4183.8
$ws.Range("N107").Value2 = -8023.8

# Row 122
$ws.Range("H122").Value2 = 1907.75
$ws.Range("I122").Value2 = 1901.6666
$ws.Range("J122").Value2 = 1999
$ws.Range("K122").Value2 = 5704.9998
$ws.Range("L122").Value2 = 5997
$ws.Range("M122").Value2 = -3254.9998
$ws.Range("N122").Value2 = -10897


# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 88
$ws.Range("H88").Value2 = 3839.05
$ws.Range("J88").Value2 = 5266.4165
$ws.Range("L88").Value2 = 5266.4165
$ws.Range("N88").Value2 = -6078.4165

# Row 91
$ws.Range("H91").Value2 = 3839.05
$ws.Range("J91").Value2 = 5266.4165
$ws.Range("L91").Value2 = 5266.4165
$ws.Range("M91").Value2 = -294
$ws.Range("N91").Value2 = -8074.4165


# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 107
$ws.Range("H107").Value2 = 1839.3846
$ws.Range("I107").Value2 = 1908.25
$ws.Range("J107").Value2 = 1013
$ws.Range("K107").Value2 = 1908.25
$ws.Range("L107").Value2 = 1013
$ws.Range("M107").Value2 = 11.75
$ws.Range("N107").Value2 = -4853


# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 132
$ws.Range("H132").Value2 = 14910
$ws.Range("I132").Value2 = 28956.5
$ws.Range("J132").Value2 = 10896.714
$ws.Range("K132").Value2 = 86869.5
$ws.Range("L132").Value2 = 32690.142
$ws.Range("M132").Value2 = -84339.5
$ws.Range("N132").Value2 = -37750.142


# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 4
$ws.Range("H4").Value2 = 7416067.5
$ws.Range("I4").Value2 = 8461957
$ws.Range("J4").Value2 = 617783.5
$ws.Range("K4").Value2 = 25385871
$ws.Range("L4").Value2 = 1853350.5
$ws.Range("M4").Value2 = -25385759
$ws.Range("N4").Value2 = -1853574.5

# Row 19
$ws.Range("H19").Value2 = 0
$ws.Range("I19").Value2 = 0
$ws.Range("K19").Value2 = 0
$ws.Range("M19").ClearContents()

# Row 80
$ws.Range("H80").Value2 = 2000
$ws.Range("J80").Value2 = 2000
$ws.Range("L80").Value2 = 6000
$ws.Range("N80").Value2 = -7872

# Row 83
$ws.Range("H83").Value2 = 2000
$ws.Range("J83").Value2 = 2000
$ws.Range("L83").Value2 = 18000
$ws.Range("N83").Value2 = -27360

# Row 86
$ws.Range("H86").Value2 = 474.75
$ws.Range("I86").Value2 = 483
$ws.Range("K86").Value2 = 1449
$ws.Range("M86").Value2 = -263

# Row 89
$ws.Range("H89").Value2 = 474.75
$ws.Range("I89").Value2 = 483
$ws.Range("K89").Value2 = 4347
$ws.Range("M89").Value2 = 1581

# Row 92
$ws.Range("H92").Value2 = 406.57144
$ws.Range("J92").Value2 = 184.6
$ws.Range("L92").Value2 = 553.8
$ws.Range("N92").Value2 = -3049.8

# Row 96
$ws.Range("H96").Value2 = 12874.75
$ws.Range("J96").Value2 = 12874.75
$ws.Range("L96").Value2 = 38624.25
$ws.Range("N96").Value2 = -42742.25

# Row 134
$ws.Range("H134").Value2 = 8027.16
$ws.Range("I134").Value2 = 7105.6523
$ws.Range("K134").Value2 = 21316.9569
$ws.Range("M134").Value2 = -16246.9569


# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 97
$ws.Range("H97").Value2 = 1419.825
$ws.Range("I97").Value2 = 1003.3125
$ws.Range("J97").Value2 = 3085.875
$ws.Range("K97").Value2 = 1003.3125
$ws.Range("L97").Value2 = 3085.875
$ws.Range("M97").Value2 = -507.3125
$ws.Range("N97").Value2 = -4077.875

# Row 122
$ws.Range("H122").Value2 = 5465
$ws.Range("I122").Value2 = 2513.5715
$ws.Range("J122").Value2 = 7760.5557
$ws.Range("K122").Value2 = 7540.7145
$ws.Range("L122").Value2 = 23281.6671
$ws.Range("M122").Value2 = -5090.7145
$ws.Range("N122").Value2 = -28181.6671


# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 45
$ws.Range("H45").Value2 = 20046
$ws.Range("I45").Value2 = 0
$ws.Range("J45").Value2 = 20046
$ws.Range("K45").Value2 = 0
$ws.Range("L45").Value2 = 20046
$ws.Range("N45").Value2 = -20860
$ws.Range("M45").ClearContents()

# Row 123
$ws.Range("H123").Value2 = 99749.5
$ws.Range("I123").Value2 = 0
$ws.Range("J123").Value2 = 99749.5
$ws.Range("K123").Value2 = 0
$ws.Range("L123").Value2 = 99749.5
$ws.Range("N123").Value2 = -109549.5
$ws.Range("M123").ClearContents()

# Row 132
$ws.Range("H132").Value2 = 2805.2424
$ws.Range("I132").Value2 = 2459.077
$ws.Range("J132").Value2 = 4091
$ws.Range("K132").Value2 = 7377.231000000001
$ws.Range("L132").Value2 = 12273
$ws.Range("M132").Value2 = -4847.231000000001
$ws.Range("N132").Value2 = -17333

# Row 141
$ws.Range("H141").Value2 = 0
$ws.Range("I141").Value2 = 0
$ws.Range("J141").Value2 = 0
$ws.Range("K141").Value2 = 0
$ws.Range("L141").Value2 = 0
$ws.Range("M141:N141").ClearContents()


# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 81
$ws.Range("H81").Value2 = 1543.5625
$ws.Range("I81").Value2 = 1257.4286
$ws.Range("J81").Value2 = 1766.1111
$ws.Range("K81").Value2 = 2514.8572
$ws.Range("L81").Value2 = 3532.2222
$ws.Range("M81").Value2 = -1453.8572
$ws.Range("N81").Value2 = -5654.2222

# Row 84
$ws.Range("H84").Value2 = 1543.5625
$ws.Range("I84").Value2 = 1257.4286
$ws.Range("J84").Value2 = 1766.1111
$ws.Range("K84").Value2 = 12574.286
$ws.Range("L84").Value2 = 17661.111
$ws.Range("M84").Value2 = -7270.286
$ws.Range("N84").Value2 = -28269.111

# Row 93
$ws.Range("H93").Value2 = 30499.75
$ws.Range("J93").Value2 = 30499.75
$ws.Range("L93").Value2 = 30499.75
$ws.Range("N93").Value2 = -35491.75

# Row 132
$ws.Range("H132").Value2 = 4935.0723
$ws.Range("I132").Value2 = 2902.0984
$ws.Range("J132").Value2 = 20436.5
$ws.Range("K132").Value2 = 8706.2952
$ws.Range("L132").Value2 = 61309.5
$ws.Range("M132").Value2 = -6176.2952
$ws.Range("N132").Value2 = -66369.5

